# Update loading_percent values on the active sheet (Sheet1)
# Diff corresponds to rerun of the "Case_1_64" (380 kV) load-flow case:
# columns B,D,E,F,G,K,L,N,O for rows 2-25 receive new loading-percent results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "B" = 16.2715326603696; "D" = 4.991257422944105; "E" = 17.67479000615848; "F" = 25.0355028630226; "G" = 3.634170068894004; "K" = 9.458621779711443; "L" = 8.4597053410284; "N" = 19.90530267709983; "O" = 22.37278902369133 }
    3 = @{ "B" = 16.16594951984914; "D" = 4.945236309674896; "E" = 17.73499294976907; "F" = 25.02023450545102; "G" = 3.636013083655973; "K" = 9.120290544060952; "L" = 8.418417321417381; "N" = 19.96866362212512; "O" = 22.41213720038299 }
    4 = @{ "B" = 16.10419250070908; "D" = 4.916288491402191; "E" = 17.77445857534432; "F" = 25.01744736491307; "G" = 3.637205626562548; "K" = 8.904087863069487; "L" = 8.394479033361748; "N" = 20.00939563577742; "O" = 22.44121863523963 }
    5 = @{ "B" = 16.07982128288002; "D" = 4.904323032505046; "E" = 17.79117064318778; "F" = 25.01797019802123; "G" = 3.637706965023858; "K" = 8.813932072025093; "L" = 8.385086867802446; "N" = 20.02645532590021; "O" = 22.45430503686284 }
    6 = @{ "B" = 16.07582313792962; "D" = 4.902326109013296; "E" = 17.79398370962946; "F" = 25.01815722658903; "G" = 3.637791141456275; "K" = 8.798840310610085; "L" = 8.383549437016226; "N" = 20.02931596170204; "O" = 22.45655257700416 }
    7 = @{ "B" = 16.10386057258916; "D" = 4.916127799031995; "E" = 17.77468141039322; "F" = 25.01744769854901; "G" = 3.637212325502556; "K" = 8.902880187816907; "L" = 8.3943508882945; "N" = 20.00962384008891; "O" = 22.44139012360639 }
    8 = @{ "B" = 16.23450424672002; "D" = 4.975535476371388; "E" = 17.69502934503031; "F" = 25.0288726819594; "G" = 3.634792924226026; "K" = 9.343765914687845; "L" = 8.445180744279556; "N" = 19.92677092574272; "O" = 22.38533366836499 }
    9 = @{ "B" = 16.51400222008566; "D" = 5.086353273881266; "E" = 17.5586443644197; "F" = 25.10340164137347; "G" = 3.630529722358969; "K" = 10.13821613752507; "L" = 8.555711034158703; "N" = 19.77873879100894; "O" = 22.31453732605181 }
    10 = @{ "B" = 16.73205804554463; "D" = 5.164048180091245; "E" = 17.47047890639757; "F" = 25.18965972164067; "G" = 3.627687885914234; "K" = 10.67576938297402; "L" = 8.643029173088911; "N" = 19.67869388366983; "O" = 22.28646959243364 }
    11 = @{ "B" = 16.83367915723084; "D" = 5.198532058117169; "E" = 17.43297518811417; "F" = 25.23565525573289; "G" = 3.626457459600263; "K" = 10.90973164349496; "L" = 8.68396420935369; "N" = 19.63505387526433; "O" = 22.27891268628292 }
    12 = @{ "B" = 16.8724799425166; "D" = 5.21146217132227; "E" = 17.41914724320026; "F" = 25.25403470515401; "G" = 3.626000444596659; "K" = 10.99676747895115; "L" = 8.699629537867057; "N" = 19.61879615624292; "O" = 22.27680073460209 }
    13 = @{ "B" = 16.86410978818048; "D" = 5.208683211824793; "E" = 17.42210871607957; "F" = 25.25003374617717; "G" = 3.626098474892662; "K" = 10.97809274232361; "L" = 8.696248598088941; "N" = 19.62228565570443; "O" = 22.27722223817907 }
    14 = @{ "B" = 16.83686509899799; "D" = 5.199598420780416; "E" = 17.4318300637129; "F" = 25.23714812506362; "G" = 3.62641968216481; "K" = 10.91692362077189; "L" = 8.685249776627092; "N" = 19.63371098273621; "O" = 22.27872390962902 }
    15 = @{ "B" = 16.82021759966294; "D" = 5.194016898492948; "E" = 17.43783334985601; "F" = 25.22938027554598; "G" = 3.62661759134886; "K" = 10.87925143646912; "L" = 8.678533741085499; "N" = 19.64074416521755; "O" = 22.27974136005685 }
    16 = @{ "B" = 16.7254631530434; "D" = 5.161776883913989; "E" = 17.47298220780561; "F" = 25.18678896773831; "G" = 3.627769547860466; "K" = 10.66026307258217; "L" = 8.640377510265361; "N" = 19.68158340307117; "O" = 22.28706835389373 }
    17 = @{ "B" = 16.66793476382221; "D" = 5.141775172435191; "E" = 17.49521131815152; "F" = 25.16238460897182; "G" = 3.628492171342324; "K" = 10.5231823033656; "L" = 8.617273417851541; "N" = 19.70711526101677; "O" = 22.29289830685654 }
    18 = @{ "B" = 16.63507616432675; "D" = 5.130190040052024; "E" = 17.50824196431438; "F" = 25.14898446096166; "G" = 3.628913675512652; "K" = 10.44334386757749; "L" = 8.604099423030576; "N" = 19.72197668327725; "O" = 22.29674205597196 }
    19 = @{ "B" = 16.62399124670378; "D" = 5.126253796879339; "E" = 17.51269601828321; "F" = 25.14455700808079; "G" = 3.62905739913809; "K" = 10.41614267879441; "L" = 8.599658975628852; "N" = 19.72703880107312; "O" = 22.29812771339024 }
    20 = @{ "B" = 16.67403514797693; "D" = 5.143912771705486; "E" = 17.49281963288289; "F" = 25.16491667501976; "G" = 3.628414639606435; "K" = 10.53787788686194; "L" = 8.619721073948927; "N" = 19.704379126867; "O" = 22.29222692841238 }
    21 = @{ "B" = 16.84485910510873; "D" = 5.20227035844109; "E" = 17.42896452415197; "F" = 25.24090692060465; "G" = 3.626325094041522; "K" = 10.93493312083173; "L" = 8.688476027084924; "N" = 19.63034782883658; "O" = 22.27826248565367 }
    22 = @{ "B" = 16.95834811891557; "D" = 5.239661255159866; "E" = 17.38941069808666; "F" = 25.29617248802397; "G" = 3.625011432142359; "K" = 11.18532053927004; "L" = 8.734362977821648; "N" = 19.58352450545446; "O" = 22.27350550190181 }
    23 = @{ "B" = 16.89761798653194; "D" = 5.219775041274268; "E" = 17.41032206486065; "F" = 25.26616718086131; "G" = 3.625707816635494; "K" = 11.05252957303963; "L" = 8.709788674290978; "N" = 19.60837262065577; "O" = 22.27564458016579 }
    24 = @{ "B" = 16.6712764915887; "D" = 5.142946630231555; "E" = 17.49390013244492; "F" = 25.16376996407421; "G" = 3.628449672829621; "K" = 10.53123721522054; "L" = 8.61861414877138; "N" = 19.70561556354499; "O" = 22.29252892597372 }
    25 = @{ "B" = 16.4360465295375; "D" = 5.057009458109247; "E" = 17.59342341278326; "F" = 25.07768243800732; "G" = 3.631631826984241; "K" = 9.931181896027484; "L" = 8.524698305738411; "N" = 19.81724862009727; "O" = 22.31453732605181 }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
